$wb = $excel.ActiveWorkbook

# --- Input sheet ---
$inputSheet = $wb.Worksheets.Item("Input")

# M7: corrected date-ish number 27/01/2022 -> 27/02/2022
$inputSheet.Range("M7").Value = 27022022

# Move the active selection to L13 (was O12)
$inputSheet.Activate()
$inputSheet.Range("L13").Select()

# --- Output sheet ---
$output = $wb.Worksheets.Item("Output")

# New timestamps recorded for this test run
$output.Range("A2").Value = "28/01/2022 2:23:14 pm"
$output.Range("A3").Value = "28/01/2022 2:23:32 pm"
$output.Range("D3").Value = "₹27,585"
$output.Range("A4").Value = "28/01/2022 2:24:02 pm"

# D4/E4/G4 must stay plain TEXT (shared string, default cell style) even
# though the values look like grouped numbers (e.g. "830,640"). A direct
# .Value assignment gets auto-parsed into a real number with a new
# comma-grouped number format/style, which does not match here, so each
# value is round-tripped through a text formula + paste-values (values
# only) instead, using a scratch cell well outside the used range.
$scratch = $output.Range("ZZ1")

$scratch.Formula = '="830,640"'
$scratch.Copy()
$output.Range("D4").PasteSpecial(-4163)

$scratch.Formula = '="38,700"'
$scratch.Copy()
$output.Range("E4").PasteSpecial(-4163)

$scratch.Formula = '="869,350"'
$scratch.Copy()
$output.Range("G4").PasteSpecial(-4163)

$scratch.Clear()

$output.Range("A5").Value = "28/01/2022 2:24:37 pm"
$output.Range("A6").Value = "28/01/2022 2:24:51 pm"
$output.Range("A8").Value = "28/01/2022 2:25:21 pm"
$output.Range("A9").Value = "28/01/2022 2:25:27 pm"
$output.Range("A10").Value = "28/01/2022 2:25:43 pm"
$output.Range("A11").Value = "28/01/2022 2:25:48 pm"
